# Federated IDM overview updates.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# --- Workbook-level changes -------------------------------------------------

# Rename the "SSO" sheet to "SSOProtocols".
$ws1.Name = "SSOProtocols"

# The workbook no longer references the external workbook
# "/Mokslai/Magistras/Islaidos.xlsx" - break/remove that external link.
[void]$wb.BreakLink("/Mokslai/Magistras/Islaidos.xlsx", 1)

# --- Sheet1 ("SSOProtocols") content updates --------------------------------

# Row 2 header: the "Savybė" label cell is now blank.
$ws1.Range("B2").Value = ""

# Row 5: "Duomenų perdavimas" row values change.
$ws1.Range("C5").Value = "HTTP, SOAP"
$ws1.Range("D5").Value = "HTTP, REST"
$ws1.Range("E5").Value = "HTTP, REST"

# Row 7: was "Pasirašymas" (signing) with empty data cells, now
# "Duomenų šifravimas" (data encryption) with all three protocols supporting it.
$ws1.Range("B7").Value = "Duomenų šifravimas"
$ws1.Range("C7").Value = "Yra"
$ws1.Range("D7").Value = "Yra"
$ws1.Range("E7").Value = "Yra"

# Row 8: was "Šifravimas" (encryption) with empty data cells, now
# "Tapatybės tiekėjo suteiktų duomenų validavimas" (IdP-supplied data validation).
$ws1.Range("B8").Value = "Tapatybės tiekėjo suteiktų duomenų validavimas"
$ws1.Range("C8").Value = "Viešo-privataus rakto infrastruktūra"
$ws1.Range("D8").Value = "Neapibrėžta (palikta realizacijai)"
$ws1.Range("E8").Value = "Viešo-privataus rakto infrastruktūra"
$ws1.Rows.Item(8).RowHeight = 43.2

# Row 11: was "Saugumo rizikos" (security risks) with empty data cells, now
# "Naudojančios organizacijos" (organizations using the protocol).
$ws1.Range("B11").Value = "Naudojančios organizacijos"
$ws1.Range("C11").Value = "Salesforce, PingFederate, Oracle Access Manager"
$ws1.Range("D11").Value = "Google, Amazon, GitHub"
$ws1.Range("E11").Value = "Google, Microsoft, Ping Identity"
$ws1.Rows.Item(11).RowHeight = 43.2

# Column width tweaks (C narrower, D a bit wider).
$ws1.Columns.Item(3).ColumnWidth = 20.17
$ws1.Columns.Item(4).ColumnWidth = 15.83

# Sheet1 becomes the active sheet/tab, with a new selected cell.
[void]$ws1.Activate()
[void]$ws1.Range("D17").Select()

# --- Sheet2 ("Vertinimas") ---------------------------------------------------
# No content changes; it simply stops being the active/selected tab, which
# happens automatically once sheet1 is activated above.
